$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds plain numeric-looking text (e.g. "542.98", "1.00") that Excel
# would otherwise auto-coerce to a Number on assignment. Force the whole data
# range to Text format first, write the values, then clear the style back to
# Normal so no stray style index is left on the cells (matches the original,
# unstyled inline-string cells).
$dCol = $ws.Range("D2:D51")
$dCol.NumberFormat = "@"

$ws.Range("D2").Value = "59.464.45"
$ws.Range("E2").Value = "  -0.74%  "
$ws.Range("D3").Value = "2.518.29"
$ws.Range("E3").Value = "  -0.76%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "542.98"
$ws.Range("E5").Value = "  -0.20%  "
$ws.Range("D6").Value = "139.89"
$ws.Range("E6").Value = "  -3.73%  "
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  +0.35%  "
$ws.Range("E8").Value = "  -1.29%  "
$ws.Range("D9").Value = "2.522.17"
$ws.Range("E9").Value = "  -1.79%  "
$ws.Range("E10").Value = "  +0.44%  "
$ws.Range("E11").Value = "  -0.07%  "
$ws.Range("D12").Value = "5.41"
$ws.Range("E12").Value = "  -2.98%  "
$ws.Range("D13").Value = "0.350"
$ws.Range("E13").Value = "  -3.26%  "
$ws.Range("D14").Value = "2.968.87"
$ws.Range("E14").Value = "  -0.62%  "
$ws.Range("D15").Value = "23.34"
$ws.Range("E15").Value = "  -1.79%  "
$ws.Range("D16").Value = "59.353.06"
$ws.Range("E16").Value = "  -0.82%  "
$ws.Range("D17").Value = "0.0000141"
$ws.Range("E17").Value = "  -1.26%  "
$ws.Range("D18").Value = "2.517.93"
$ws.Range("E18").Value = "  -1.57%  "
$ws.Range("D19").Value = "11.10"
$ws.Range("E19").Value = "  -1.69%  "
$ws.Range("D20").Value = "4.29"
$ws.Range("E20").Value = "  -0.78%  "
$ws.Range("D21").Value = "325.45"
$ws.Range("E21").Value = "  -0.88%  "
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("D23").Value = "5.87"
$ws.Range("E23").Value = "  -1.29%  "
$ws.Range("D24").Value = "63.45"
$ws.Range("E24").Value = "  +1.88%  "
$ws.Range("D25").Value = "0.421"
$ws.Range("E25").Value = "  -4.66%  "
$ws.Range("D26").Value = "0.168"
$ws.Range("E26").Value = "  +1.21%  "
$ws.Range("D27").Value = "1.00"
$ws.Range("E27").Value = "  +0.96%  "
$ws.Range("D28").Value = "7.78"
$ws.Range("E28").Value = "  -3.06%  "
$ws.Range("D29").Value = "0.0₃0783"
$ws.Range("E29").Value = "  -1.84%  "
$ws.Range("D30").Value = "6.73"
$ws.Range("E30").Value = "  -4.14%  "
$ws.Range("E31").Value = "  -1.09%  "
$ws.Range("D32").Value = "164.11"
$ws.Range("E32").Value = "  +1.10%  "
$ws.Range("D34").Value = "1.12"
$ws.Range("E34").Value = "  -8.10%  "
$ws.Range("D35").Value = "1.43"
$ws.Range("E35").Value = "  -2.39%  "
$ws.Range("D36").Value = "18.56"
$ws.Range("E36").Value = "  -1.46%  "
$ws.Range("D37").Value = "4.20"
$ws.Range("E37").Value = "  -6.07%  "
$ws.Range("D38").Value = "1.60"
$ws.Range("E38").Value = "  -2.18%  "
$ws.Range("D39").Value = "3.67"
$ws.Range("E39").Value = "  -1.61%  "
$ws.Range("D40").Value = "0.809"
$ws.Range("E40").Value = "  -3.50%  "
$ws.Range("D41").Value = "5.21"
$ws.Range("E41").Value = "  -8.51%  "
$ws.Range("D42").Value = "280.93"
$ws.Range("E42").Value = "  -7.09%  "
$ws.Range("D43").Value = "0.998"
$ws.Range("E43").Value = "  +0.50%  "
$ws.Range("B44").Value = "Mantle"
$ws.Range("C44").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D44").Value = "0.600"
$ws.Range("E44").Value = "  -1.52%  "
$ws.Range("B45").Value = "WhiteBITCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D45").Value = "10.87"
$ws.Range("E45").Value = "  +0.29%  "
$ws.Range("D46").Value = "0.0939"
$ws.Range("E46").Value = "  +0.10%  "
$ws.Range("D47").Value = "124.82"
$ws.Range("E47").Value = "  +0.00%  "
$ws.Range("D48").Value = "0.0513"
$ws.Range("E48").Value = "  -1.37%  "
$ws.Range("D49").Value = "0.0224"
$ws.Range("E49").Value = "  -2.32%  "
$ws.Range("D50").Value = "17.95"
$ws.Range("E50").Value = "  -2.14%  "
$ws.Range("D51").Value = "1.773.58"
$ws.Range("E51").Value = "  -2.58%  "

$dCol.Style = "Normal"
